$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 40
$ws.Range("F5").Value = 71
$ws.Range("F6").Value = 2399
$ws.Range("F9").Value = 44
$ws.Range("F14").Value = 836
$ws.Range("F15").Value = 1092
$ws.Range("F17").Value = 3406
$ws.Range("F21").Value = 722
$ws.Range("F26").Value = 1095
$ws.Range("F29").Value = 919
$ws.Range("F30").Value = 895

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 200

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 468

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 40
$ws.Range("F11").Value = 71
$ws.Range("F12").Value = 468
$ws.Range("F13").Value = 2399
$ws.Range("F16").Value = 44
$ws.Range("F28").Value = 836
$ws.Range("F30").Value = 200
$ws.Range("F31").Value = 1092
$ws.Range("F34").Value = 3406
$ws.Range("F37").Value = 722
$ws.Range("F41").Value = 1095
$ws.Range("F49").Value = 919
$ws.Range("F50").Value = 895
